# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The curation moved several fields that used to be modelled as
# "dimensions" (sdmx-dimension / iaest-dimension, skos:Concept typed,
# with an accompanying mapping-*.xlsx lookup file) into plain
# "measures" (iaest-measure, xsd:int typed, no mapping file):
#   - residencia-continente-nombre (column C)
#   - inscripcion-provincia-nombre (column H)
#   - inscripcion-municipio-estrato (column I)
#   - residencia-area-nombre (column J)
#   - sexo (column M)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-measure / sdmx-dimension qualifier per column
$ws.Range("C2").Value = "iaest-measure:residencia-continente-nombre"
$ws.Range("H2").Value = "iaest-measure:inscripcion-provincia-nombre"
$ws.Range("I2").Value = "iaest-measure:inscripcion-municipio-estrato"
$ws.Range("J2").Value = "iaest-measure:residencia-area-nombre"
$ws.Range("M2").Value = "iaest-measure:sexo"

# Row 3: dim -> medida
$ws.Range("C3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "medida"
$ws.Range("M3").Value = "medida"

# Row 4: skos:Concept / URI-* -> xsd:int
$ws.Range("C4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("M4").Value = "xsd:int"

# Row 5: mapping files for these measures no longer apply, clear them
$ws.Range("C5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("M5").Value = ""
